# Applies the "Allowed CMS template to have second set of cpt" edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Claim sheet
# ---------------------------------------------------------------------
$claim = $wb.Worksheets.Item("Claim")

# Column A (labels) is locked by the sheet protection, so unprotect
# temporarily to update the label formulas/text, then re-protect.
$claim.Unprotect()

$claim.Range("A12").Formula = '=IF($B$2="Professional (CMS)", "Modifier", IF($B$2="Institutional (UB)", "Description (1)", "Error"))'
$claim.Range("A13").Formula = '=IF($B$2="Professional (CMS)", "Diagnosis", IF($B$2="Institutional (UB)", "CPT Code (1)", "Error"))'
$claim.Range("A14").Formula = '=IF($B$2="Professional (CMS)","CPT Code (1)",IF($B$2="Institutional (UB)", "Charges (1)", "Error"))'
$claim.Range("A15").Formula = '=IF($B$2="Professional (CMS)", "Charges (1)", IF($B$2="Institutional (UB)", "Units (1)", "Error"))'
$claim.Range("A16").Formula = '=IF($B$2="Professional (CMS)", "Units (1)", IF($B$2="Institutional (UB)", "Description (2)", "Error"))'

# Rows 17-19 switch from blank formulas to plain literal labels for the
# new "second set" of CPT/Charges/Units fields.
$claim.Range("A17").Value = "CPT Code (2)"
$claim.Range("A18").Value = "Charges (2)"
$claim.Range("A19").Value = "Units (2)"

$claim.Protect()

# Column B values (these cells are unlocked even while protected).
$claim.Range("B3").Value = ""
$claim.Range("B4").Value = ""

$claim.Range("B12").Value = "U1"
$claim.Range("B13").Value = "A"
$claim.Range("B14").Value = "S5105"
$claim.Range("B15").Value = "75"
# B16 (Units (1)) keeps its existing value of 1 - no change required.

# B19's format needs to switch from the general-number style to the
# shared "centered text" style already used by B17 - this only requires
# a number-format change, which lets the engine reuse the existing style.
$claim.Range("B19").NumberFormat = "@"

# ---------------------------------------------------------------------
# Members sheet
# ---------------------------------------------------------------------
$members = $wb.Worksheets.Item("Members")

$members.Range("I2").Value = 45747
$members.Range("I3").Value = 45731

$members.Range("I4").Select() | Out-Null
